# "Generate Report for Handoff"
# Regenerates the localization-status handoff report: for every source file
# that is currently "Ready for handoff" (rows 5-16 on both the zh-cn and
# de-de sheets), the report's Priority flips from "ht" to "mt" and the
# Latest Handoff Datetime is bumped to the new generation timestamp.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
for ($r = 5; $r -le 16; $r++) {
    $ws_overview.Cells.Item($r, 7).Value = "2016-08-20 16:33:56"
}

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 5; $r -le 16; $r++) {
    $ws_zhcn.Cells.Item($r, 5).Value = "mt"
    $ws_zhcn.Cells.Item($r, 8).Value = "2016-08-20 16:33:51"
}

$ws_dede = $wb.Worksheets.Item("de-de")
for ($r = 5; $r -le 16; $r++) {
    $ws_dede.Cells.Item($r, 5).Value = "mt"
    $ws_dede.Cells.Item($r, 8).Value = "2016-08-20 16:33:56"
}
